$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @{ D='27.389.48'; E='  -2.42%  ' }
    3 = @{ D='1.651.25'; E='  -2.28%  ' }
    4 = @{ E='  -0.14%  ' }
    5 = @{ D='''213.66'; E='  -1.41%  ' }
    6 = @{ E='  -1.93%  ' }
    7 = @{ E='  -0.09%  ' }
    8 = @{ D='''24.20'; E='  +0.23%  ' }
    10 = @{ D='''0.0615'; E='  -1.69%  ' }
    11 = @{ D='''0.0877'; E='  -0.89%  ' }
    12 = @{ B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='1.885.12'; E='  -2.32%  ' }
    13 = @{ B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.645.61'; E='  -2.61%  ' }
    14 = @{ D='''0.575'; E='  +3.07%  ' }
    15 = @{ E='  -2.34%  ' }
    16 = @{ D='''65.96'; E='  -1.41%  ' }
    17 = @{ D='27.395.95'; E='  -2.27%  ' }
    18 = @{ D='''234.15'; E='  -6.55%  ' }
    19 = @{ D='0.0₃0727'; E='  -2.20%  ' }
    20 = @{ D='''7.46'; E='  -2.99%  ' }
    21 = @{ E='  -0.17%  ' }
    22 = @{ D='''4.39'; E='  -3.16%  ' }
    23 = @{ E='  -2.49%  ' }
    24 = @{ D='''2.01'; E='  -1.59%  ' }
    25 = @{ D='''146.94'; E='  -0.52%  ' }
    26 = @{ D='''7.20'; E='  -1.84%  ' }
    27 = @{ D='''16.05'; E='  -2.73%  ' }
    28 = @{ E='  +0.00%  ' }
    29 = @{ E='  -2.15%  ' }
    30 = @{ E='  -1.29%  ' }
    31 = @{ E='  -4.18%  ' }
    32 = @{ E='  -2.19%  ' }
    33 = @{ D='1.458.85'; E='  +0.83%  ' }
    34 = @{ D='''3.10'; E='  -2.65%  ' }
    35 = @{ E='  -3.99%  ' }
    36 = @{ D='''2.38'; E='  -0.78%  ' }
    37 = @{ E='  -3.98%  ' }
    38 = @{ E='  -3.35%  ' }
    39 = @{ E='  -1.42%  ' }
    40 = @{ E='  -0.17%  ' }
    41 = @{ E='  -0.12%  ' }
    42 = @{ D='''65.51' }
    43 = @{ D='''5.44'; E='  -1.68%  ' }
    44 = @{ D='''2.22'; E='  -0.86%  ' }
    45 = @{ D='1.793.80'; E='  -2.34%  ' }
    46 = @{ E='  -1.71%  ' }
    47 = @{ E='  -0.31%  ' }
    48 = @{ D='''88.59'; E='  -1.06%  ' }
    49 = @{ E='  -4.56%  ' }
    50 = @{ E='  -1.81%  ' }
    51 = @{ D='''7.78'; E='  -1.96%  ' }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    foreach ($col in $vals.Keys) {
        $addr = "$col$r"
        $ws.Range($addr).Value = $vals[$col]
    }
}

Write-Host "Applied $($rowData.Count) row updates"